$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.699.09'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.879.34'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4756'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2824'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06492'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +15.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.869.37'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07578'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.13'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +13.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.037'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6475'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '297.70'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +31.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.682.53'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007489'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.136.45'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.121'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.105'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.27'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.208'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.51'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +9.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.943'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.28%  '
$ws.Range("E29").Value = '  +2.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.349'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.141'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.937'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04967'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.165'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7161'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.718'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01899'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.698'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.043'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +6.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8955'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '106.80'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  +4.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.565'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.31'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +8.58%  '
$ws.Range("E46").Value = '  +4.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1212'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '34.43'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.792'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05610'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.379'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.11%  '
